$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes ---------------------------------------------------

# Row 9: "Entity"/"prov" pair (G9/I9) becomes "Event"/"sem"
$ws.Range("G9").Value = "Event"
$ws.Range("I9").Value = "sem"

# Row 10/11: K column "leolani talk" -> "leolani world"
$ws.Range("K10").Value = "leolani world"
$ws.Range("K11").Value = "leolani world"

# New cell L10: add note about Event class
$ws.Range("L10").Value = "has time and has actor and has subevent"

# --- Style changes -------------------------------------------------------

# D4, D5, D10 ("Present?" y cells) lose their bold + left/center alignment,
# reverting to plain/default formatting.
$ws.Range("D4").ClearFormats()
$ws.Range("D4").Value = "y"

$ws.Range("D5").ClearFormats()
$ws.Range("D5").Value = "y"

$ws.Range("D10").ClearFormats()
$ws.Range("D10").Value = "y"

# Row 13 (A13:D13) gains bold + left/center alignment (matching the style
# used elsewhere for the same kind of header-like row, e.g. row 4/5/10).
$r13 = $ws.Range("A13:D13")
$r13.Font.Bold = $true
$r13.HorizontalAlignment = -4131
$r13.VerticalAlignment = -4108

# --- Selection -------------------------------------------------------

$ws.Range("K11").Select()
